# Apply weekly timesheet corrections for Chris Jacobi's 2026-01-19 export.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Weekly Timesheet")

# --- Client name corrections (column B) ---
$ws.Range("B2").Value = "Knight"
$ws.Range("B3").Value = "McClure"
$ws.Range("B4").Value = "McGill"
$ws.Range("B5").Value = "TOTAL:"
$ws.Range("B6").Value = "Moulton"

# --- Hours / Rate / Total corrections ---
$ws.Range("C2").Value = 9.5
$ws.Range("F2").Value = 950

$ws.Range("C3").Value = 7
$ws.Range("F3").Value = 700

$ws.Range("E4").Value = 90
$ws.Range("F4").Value = 630

$ws.Range("C5").Value = 6.5
$ws.Range("F5").Value = 650

$ws.Range("C6").Value = 8
$ws.Range("F6").Value = 800

# --- Updated subtotal total ---
$ws.Range("F8").Value = 3730
